$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Available"
$ws.Range("B2").Value = "SPA"
$ws.Range("C2").Value = "WB"
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 152
$ws.Range("G2").Value = 43640.51221876157
$ws.Range("G2").NumberFormat = "m/d/yy h:mm"

$ws.Range("A3").Value = "Available"
$ws.Range("B3").Value = "SPA"
$ws.Range("C3").Value = "WB"
$ws.Range("D3").Value = 15
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 589
$ws.Range("G3").Value = 43640.51225010417
$ws.Range("G3").NumberFormat = "m/d/yy h:mm"
